$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.068.71"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "2.301.74"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.76"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.70"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.29"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.69%  "

$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +12.97%  "

$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").Value = "2.660.81"
$ws.Range("E16").Value = "  +1.38%  "

$ws.Range("D17").Value = "2.308.23"
$ws.Range("E17").Value = "  +1.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.813"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.65%  "

$ws.Range("D19").Value = "43.023.57"
$ws.Range("E19").Value = "  +1.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.75"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.82%  "

$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("E22").Value = "  +1.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.22%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -1.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.44"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +10.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.34"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.12"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("E35").Value = "  +6.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.73"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.07%  "

$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.83"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("E42").Value = "  +0.19%  "

$ws.Range("E43").Value = "  -0.90%  "

$ws.Range("D44").Value = "1.980.46"
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.94%  "

$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("D49").Value = "2.528.52"
$ws.Range("E49").Value = "  +1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.38"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.59"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.97%  "
